# metadata.xlsx -- "Add files via upload" (17/06/25)
#
# On the "sample_annotation" sheet:
#  - column A (rows 2-9) sample ids drop their "_PC346C" / "_PC339" suffix,
#    e.g. "001_PC346C" -> "001", "005_PC339" -> "005"
#  - F1 loses its (redundant, no-op) explicit fill/border formatting, keeping
#    just the bold font it already shares with the other header cells
#  - the AutoFilter on the sheet is removed
#  - the sheet's remembered selection moves from E1 to F9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sample_annotation")

$ws.Range("A2").Value = "'001"
$ws.Range("A3").Value = "'002"
$ws.Range("A4").Value = "'003"
$ws.Range("A5").Value = "'004"
$ws.Range("A6").Value = "'005"
$ws.Range("A7").Value = "'006"
$ws.Range("A8").Value = "'007"
$ws.Range("A9").Value = "'008"

# Drop F1's explicit (but visually no-op) fill/border formatting.
$ws.Range("F1").Interior.Pattern = -4142  # xlNone
$ws.Range("F1").Borders.LineStyle = -4142 # xlNone

# Remove the sheet's AutoFilter entirely.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# Leave the sheet with F9 as the remembered selection.
$ws.Range("F9").Select()

$wb.Save()
